$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily ventilator records to append (Estado, Fecha, Cantidad)
$newRows = @(
    @("disponibles", "2022-09-14", 293),
    @("ocupados",    "2022-09-14", 1649),
    @("disponibles", "2022-09-15", 287),
    @("ocupados",    "2022-09-15", 1652),
    @("disponibles", "2022-09-16", 296),
    @("ocupados",    "2022-09-16", 1643),
    @("disponibles", "2022-09-17", 315),
    @("ocupados",    "2022-09-17", 1628),
    @("disponibles", "2022-09-18", 309),
    @("ocupados",    "2022-09-18", 1638),
    @("disponibles", "2022-09-19", 341),
    @("ocupados",    "2022-09-19", 1607),
    @("disponibles", "2022-09-20", 303),
    @("ocupados",    "2022-09-20", 1644),
    @("disponibles", "2022-09-21", 275),
    @("ocupados",    "2022-09-21", 1660),
    @("disponibles", "2022-09-22", 276),
    @("ocupados",    "2022-09-22", 1660)
)

$startRow = 1768
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $estado = $newRows[$i][0]
    $fecha = $newRows[$i][1]
    $cantidad = $newRows[$i][2]

    $ws.Range("A$row").Value = $estado

    # Column B stores dates as plain text (e.g. "2022-09-14"), not real
    # dates, matching the rest of the sheet. Force text so Excel doesn't
    # auto-convert the string into a date serial number, then clear the
    # temporary formatting so the cell keeps the sheet's default style.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $fecha
    $ws.Range("B$row").ClearFormats()

    $ws.Range("C$row").Value = $cantidad
}
